$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("FeltMotion_Prob")

# Sheet1 (raw trial data) corrections
$ws1.Range("G2").Value = 1
$ws1.Range("F4").Value = -1
$ws1.Range("I4").Value = 1
$ws1.Range("F6").Value = -1
$ws1.Range("I6").Value = 1
$ws1.Range("G7").Value = 1
$ws1.Range("F8").Value = -1
$ws1.Range("I8").Value = 1
$ws1.Range("F9").Value = -1
$ws1.Range("I9").Value = 1
$ws1.Range("F11").Value = 0
$ws1.Range("G11").Value = 1
$ws1.Range("I11").Value = 1
$ws1.Range("F12").Value = -1
$ws1.Range("I12").Value = 1
$ws1.Range("F13").Value = 0
$ws1.Range("G13").Value = 1
$ws1.Range("I13").Value = 1
$ws1.Range("F14").Value = 0
$ws1.Range("I14").Value = 1
$ws1.Range("F15").Value = 1
$ws1.Range("I15").Value = 1
$ws1.Range("F16").Value = 1
$ws1.Range("I16").Value = 1
$ws1.Range("F17").Value = 0
$ws1.Range("I17").Value = 1
$ws1.Range("F18").Value = 0
$ws1.Range("I18").Value = 1
$ws1.Range("F19").Value = 1
$ws1.Range("I19").Value = 1
$ws1.Range("G20").Value = 1
$ws1.Range("F21").Value = 0
$ws1.Range("I21").Value = 1
$ws1.Range("F22").Value = 1
$ws1.Range("I22").Value = 1
$ws1.Range("F23").Value = 1
$ws1.Range("I23").Value = 1
$ws1.Range("F24").Value = 1
$ws1.Range("I24").Value = 1
$ws1.Range("G25").Value = 1
$ws1.Range("F26").Value = 1
$ws1.Range("I26").Value = 1
$ws1.Range("F28").Value = 1
$ws1.Range("I28").Value = 1
$ws1.Range("F29").Value = 0
$ws1.Range("I29").Value = 0
$ws1.Range("F30").Value = 1
$ws1.Range("I30").Value = 1
$ws1.Range("F31").Value = 0
$ws1.Range("I31").Value = 0
$ws1.Range("F32").Value = -1
$ws1.Range("I32").Value = 1
$ws1.Range("F33").Value = -1
$ws1.Range("I33").Value = 1
$ws1.Range("F34").Value = -1
$ws1.Range("I34").Value = 1
$ws1.Range("F35").Value = -1
$ws1.Range("I35").Value = 1
$ws1.Range("F38").Value = 0
$ws1.Range("I38").Value = 1
$ws1.Range("G53").Value = 1
$ws1.Range("G73").Value = 1
$ws1.Range("G78").Value = 1
$ws1.Range("G80").Value = 1
$ws1.Range("G90").Value = 1
$ws1.Range("G98").Value = 1
$ws1.Range("G121").Value = 1
$ws1.Range("G123").Value = 1
$ws1.Range("G134").Value = 1
$ws1.Range("G137").Value = 1
$ws1.Range("G146").Value = 1
$ws1.Range("G181").Value = 1

# FeltMotion_Prob (summary) corrections
$ws2.Range("D2").Value = 0.1111111111111111
$ws2.Range("E2").Value = 0.1111111111111111
$ws2.Range("D3").Value = 0.7777777777777778
$ws2.Range("E3").Value = 0.1469861839480328
$ws2.Range("D5").Value = 0.75
$ws2.Range("E5").Value = 0.1636634176769943
$ws2.Range("D6").Value = 0.6666666666666666
$ws2.Range("E6").Value = 0.1666666666666667
$ws2.Range("D7").Value = 0.875
$ws2.Range("E7").Value = 0.125
$ws2.Range("D14").Value = 1
$ws2.Range("D18").Value = 0.9
$ws2.Range("E18").Value = 0.09999999999999999
